$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 84

# Assign the date as literal text (leading apostrophe prevents Excel from
# auto-converting it to a date serial number), then clear any formatting
# residue (e.g. quote-prefix style) left behind so the cell matches the
# plain, unstyled data cells above it.
$ws.Cells.Item($row, 1).Formula = "'11/24/2025"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 0.2067077544858842
$ws.Cells.Item($row, 3).Value = 0.7932922455141158
